$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 79, pushing the existing rows 79:153 down to 80:154
$ws.Rows.Item(79).Insert()

# Populate the new row 79 with the new weekly price record
$ws.Cells.Item(79, 1).Value = 4
$ws.Cells.Item(79, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(79, 3).Value = "Los Lagos"
$ws.Cells.Item(79, 4).Value = 44484
$ws.Cells.Item(79, 5).Value = 10
$ws.Cells.Item(79, 6).Value = "Fruta"
$ws.Cells.Item(79, 7).Value = 100102
$ws.Cells.Item(79, 8).Value = "Cítricos"
$ws.Cells.Item(79, 9).Value = 100102006
$ws.Cells.Item(79, 10).Value = "Pomelo"
$ws.Cells.Item(79, 11).Value = "Start Ruby"
$ws.Cells.Item(79, 12).Value = "Primera"
$ws.Cells.Item(79, 13).Value = 240
$ws.Cells.Item(79, 14).Value = 11000
$ws.Cells.Item(79, 15).Value = 12000
$ws.Cells.Item(79, 16).Value = 11500
$ws.Cells.Item(79, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(79, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(79, 19).Value = 821
$ws.Cells.Item(79, 20).Value = 14
